# Account Owner Drill Down
# Updates Monthly Billing / Monthly Collection / Quarterly ARR & Service Rev figures,
# adds a 0.00 number format to the (currently empty) "D" helper column on the
# Monthly Billing sheet, and moves the active sheet/selection from
# "Account Owners" to "Monthly Collection".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Monthly Billing ("Achievement INR Cr" column C, plus new D helper column)
# ---------------------------------------------------------------------------
$wsBilling = $wb.Worksheets.Item("Monthly Billing")

$wsBilling.Range("C4").Value  = 4.233365
$wsBilling.Range("C5").Value  = 14.377390500000001
$wsBilling.Range("C6").Value  = 13.1970253
$wsBilling.Range("C7").Value  = 20.910151599999999
$wsBilling.Range("C8").Value  = 5.5789847000000004
$wsBilling.Range("C9").Value  = 21.392120899999998
$wsBilling.Range("C10").Value = 8.1453403000000009
$wsBilling.Range("C11").Value = 67.188736399999996
$wsBilling.Range("C12").Value = 24.515263399999998
$wsBilling.Range("C13").Value = 14.112248599999999
$wsBilling.Range("C14").Value = 6.4880142999999997

# D4:D14 were blank style-1 cells; give them a 0.00 number format (new cellXfs entry)
$wsBilling.Range("D4:D14").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Monthly Collection ("Achievement INR Cr" column C)
# ---------------------------------------------------------------------------
$wsCollection = $wb.Worksheets.Item("Monthly Collection")

$wsCollection.Range("C4").Value  = 10.1832961
$wsCollection.Range("C5").Value  = 15.8174274
$wsCollection.Range("C6").Value  = 11.746712
$wsCollection.Range("C7").Value  = 14.619760100000001
$wsCollection.Range("C8").Value  = 10.1704145
$wsCollection.Range("C9").Value  = 23.058410299999998
$wsCollection.Range("C10").Value = 14.072506799999999
$wsCollection.Range("C11").Value = 38.941495600000003
$wsCollection.Range("C12").Value = 67.316375899999997
$wsCollection.Range("C13").Value = 21.3923387
$wsCollection.Range("C14").Value = 2.3051662999999998

# ---------------------------------------------------------------------------
# Quarterly ARR & Service Rev ("Service Rev Achievement" column E)
# ---------------------------------------------------------------------------
$wsArr = $wb.Worksheets.Item("Quarterly ARR & Service Rev")

$wsArr.Range("E4").Value = 27.614460697000002
$wsArr.Range("E5").Value = 29.597379536999998
$wsArr.Range("E6").Value = 24.333428848000011
$wsArr.Range("E7").Value = 22.434157786999993

# ---------------------------------------------------------------------------
# Selections on sheets that do not become the active tab.
# (Selecting a range on a sheet activates it, so these are done before the
# final activation of "Monthly Collection" below.)
# ---------------------------------------------------------------------------
$wsBilling.Range("E4").Select() | Out-Null
$wsArr.Range("E5").Select() | Out-Null

$wsOwners = $wb.Worksheets.Item("Account Owners")
$wsOwners.Range("D4").Select() | Out-Null

# ---------------------------------------------------------------------------
# Finally activate "Monthly Collection" and select B4 - this becomes the
# workbook's active tab / tabSelected sheet, replacing "Account Owners".
# ---------------------------------------------------------------------------
$wsCollection.Activate() | Out-Null
$wsCollection.Range("B4").Select() | Out-Null
